# Preprocess: fix row 4 values that were shifted (swap D4<->H4 and E4<->G4)
# add some comment to preprocess
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "C9"
$ws.Range("E4").Value = "C8"
$ws.Range("G4").Value = "C6"
$ws.Range("H4").Value = "C5"
